# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) across several Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1838.5
$ws.Range("I88").Value = 809.3333
$ws.Range("J88").Value = 2224.4375
$ws.Range("K88").Value = 809.3333
$ws.Range("L88").Value = 2224.4375
$ws.Range("M88").Value = -403.3333
$ws.Range("N88").Value = -3036.4375

$ws.Range("H91").Value = 1838.5
$ws.Range("I91").Value = 809.3333
$ws.Range("J91").Value = 2224.4375
$ws.Range("K91").Value = 809.3333
$ws.Range("L91").Value = 2224.4375
$ws.Range("M91").Value = 594.6667
$ws.Range("N91").Value = -5032.4375

$ws.Range("H112").Value = 1178589.6
$ws.Range("J112").Value = 1541090.4
$ws.Range("L112").Value = 4623271.199999999
$ws.Range("N112").Value = -4625487.199999999

$ws.Range("H129").Value = 1334.5
$ws.Range("I129").Value = 541.8
$ws.Range("J129").Value = 1664.7916
$ws.Range("K129").Value = 1625.4
$ws.Range("L129").Value = 4994.3748
$ws.Range("M129").Value = 3374.6
$ws.Range("N129").Value = -14994.3748

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 7464.1
$ws.Range("I28").Value = 4737.8887
$ws.Range("J28").Value = 32000
$ws.Range("K28").Value = 4737.8887
$ws.Range("L28").Value = 32000
$ws.Range("M28").Value = -4545.8887
$ws.Range("N28").Value = -32384

$ws.Range("H32").Value = 16328.912
$ws.Range("I32").Value = 3558.7703
$ws.Range("J32").Value = 173827.33
$ws.Range("K32").Value = 3558.7703
$ws.Range("L32").Value = 173827.33
$ws.Range("M32").Value = -3271.7703
$ws.Range("N32").Value = -174401.33

$ws.Range("H55").Value = 22584.75
$ws.Range("J55").Value = 22584.75
$ws.Range("L55").Value = 22584.75
$ws.Range("N55").Value = -23214.75

$ws.Range("H80").Value = 26961.75
$ws.Range("J80").Value = 26961.75
$ws.Range("L80").Value = 26961.75
$ws.Range("N80").Value = -28957.75

$ws.Range("H83").Value = 26961.75
$ws.Range("J83").Value = 26961.75
$ws.Range("L83").Value = 80885.25
$ws.Range("N83").Value = -90869.25

$ws.Range("H99").Value = 7464.1
$ws.Range("I99").Value = 4737.8887
$ws.Range("J99").Value = 32000
$ws.Range("K99").Value = 4737.8887
$ws.Range("L99").Value = 32000
$ws.Range("M99").Value = -1742.8887
$ws.Range("N99").Value = -37990

$ws.Range("H122").Value = 1442
$ws.Range("I122").Value = 1244.125
$ws.Range("J122").Value = 1758.6
$ws.Range("K122").Value = 3732.375
$ws.Range("L122").Value = 5275.799999999999
$ws.Range("M122").Value = -1282.375
$ws.Range("N122").Value = -10175.8

$ws.Range("H124").Value = 21666.666
$ws.Range("J124").Value = 21666.666
$ws.Range("L124").Value = 21666.666
$ws.Range("N124").Value = -31486.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17364.6
$ws.Range("I82").Value = 8128.5
$ws.Range("J82").Value = 23522
$ws.Range("K82").Value = 8128.5
$ws.Range("L82").Value = 23522
$ws.Range("M82").Value = -7745.5
$ws.Range("N82").Value = -24288

$ws.Range("H85").Value = 17364.6
$ws.Range("I85").Value = 8128.5
$ws.Range("J85").Value = 23522
$ws.Range("K85").Value = 8128.5
$ws.Range("L85").Value = 23522
$ws.Range("M85").Value = -6802.5
$ws.Range("N85").Value = -26174

$ws.Range("H92").Value = 35000
$ws.Range("J92").Value = 35000
$ws.Range("L92").Value = 35000
$ws.Range("N92").Value = -39992

$ws.Range("H94").Value = 843.76
$ws.Range("I94").Value = 739.1111
$ws.Range("J94").Value = 1112.8572
$ws.Range("K94").Value = 739.1111
$ws.Range("L94").Value = 1112.8572
$ws.Range("M94").Value = -288.1111
$ws.Range("N94").Value = -2014.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 13981.5
$ws.Range("J41").Value = 18232.5
$ws.Range("L41").Value = 18232.5
$ws.Range("N41").Value = -19088.5

$ws.Range("H50").Value = 8220.272000000001
$ws.Range("J50").Value = 8642.299999999999
$ws.Range("L50").Value = 8642.299999999999
$ws.Range("N50").Value = -9892.299999999999

$ws.Range("H51").Value = 7525
$ws.Range("J51").Value = 10550
$ws.Range("L51").Value = 10550
$ws.Range("N51").Value = -12022

$ws.Range("H61").Value = 7525
$ws.Range("J61").Value = 10550
$ws.Range("L61").Value = 10550
$ws.Range("N61").Value = -11246

$ws.Range("H68").Value = 35171.25
$ws.Range("J68").Value = 35171.25
$ws.Range("L68").Value = 35171.25
$ws.Range("N68").Value = -36669.25

$ws.Range("H71").Value = 35171.25
$ws.Range("J71").Value = 35171.25
$ws.Range("L71").Value = 105513.75
$ws.Range("N71").Value = -113001.75

$ws.Range("H105").Value = 2667.2222
$ws.Range("I105").Value = 1189.091
$ws.Range("K105").Value = 1189.091
$ws.Range("M105").Value = 557.9090000000001

$ws.Range("H109").Value = 20594
$ws.Range("J109").Value = 20594
$ws.Range("L109").Value = 20594
$ws.Range("N109").Value = -22674

$ws.Range("H132").Value = 3547.3914
$ws.Range("I132").Value = 3059.25
$ws.Range("K132").Value = 9177.75
$ws.Range("M132").Value = -6647.75

$ws.Range("H134").Value = 6768.0435
$ws.Range("I134").Value = 9644.786
$ws.Range("J134").Value = 2293.111
$ws.Range("K134").Value = 28934.358
$ws.Range("L134").Value = 6879.333
$ws.Range("M134").Value = -26399.358
$ws.Range("N134").Value = -11949.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 628.0303
$ws.Range("I5").Value = 361.2
$ws.Range("K5").Value = 1083.6
$ws.Range("M5").Value = -971.5999999999999

$ws.Range("H94").Value = 9314.666999999999
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 10977.6
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 32932.8
$ws.Range("M94").Value = -2324
$ws.Range("N94").Value = -34284.8

$ws.Range("H122").Value = 2657.49
$ws.Range("I122").Value = 208.38889
$ws.Range("J122").Value = 3195.0977
$ws.Range("K122").Value = 1875.50001
$ws.Range("L122").Value = 28755.8793
$ws.Range("M122").Value = 574.49999
$ws.Range("N122").Value = -33655.8793

$ws.Range("H131").Value = 842.5333000000001
$ws.Range("I131").Value = 360
$ws.Range("J131").Value = 877
$ws.Range("K131").Value = 1080
$ws.Range("L131").Value = 2631
$ws.Range("M131").Value = 3960
$ws.Range("N131").Value = -12711

$ws.Range("H135").Value = 628.0303
$ws.Range("I135").Value = 361.2
$ws.Range("K135").Value = 3250.8
$ws.Range("M135").Value = -715.7999999999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 21530.4
$ws.Range("J123").Value = 21530.4
$ws.Range("L123").Value = 21530.4
$ws.Range("N123").Value = -26430.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2495.68
$ws.Range("I132").Value = 2085.6
$ws.Range("K132").Value = 6256.799999999999
$ws.Range("M132").Value = -3726.799999999999
